# Weekly data refresh: a new daily price observation (2021-09-28) is added
# to the "Acelga" / Femacal de La Calera price series. The new record is
# inserted as row 87, pushing all existing records (old rows 87-199) down
# by one (to rows 88-200), exactly like Excel's native "Insert Row" / row
# shift-down behaviour.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 87 (shifts rows 87:199 -> 88:200,
# growing the used range from A1:R199 to A1:R200).
$ws.Rows(87).Insert()

# Populate the newly inserted row 87 with the new day's observation.
$ws.Cells.Item(87, 1).Value  = 3
$ws.Cells.Item(87, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(87, 3).Value  = "Coquimbo"
$ws.Cells.Item(87, 4).Value  = 44467
$ws.Cells.Item(87, 5).Value  = 5
$ws.Cells.Item(87, 6).Value  = 100112009
$ws.Cells.Item(87, 7).Value  = "Acelga"
$ws.Cells.Item(87, 8).Value  = "Sin especificar"
$ws.Cells.Item(87, 9).Value  = "Primera"
$ws.Cells.Item(87, 10).Value = 250
$ws.Cells.Item(87, 11).Value = 2000
$ws.Cells.Item(87, 12).Value = 2200
$ws.Cells.Item(87, 13).Value = 2104
$ws.Cells.Item(87, 14).Value = "`$/docena de atados (6 kilos)"
$ws.Cells.Item(87, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(87, 16).Value = 351
$ws.Cells.Item(87, 17).Value = 6
$ws.Cells.Item(87, 18).Value = "Hortaliza"
